$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '97.709.07'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.49%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.346.98'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -2.41%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '250.27'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -3.48%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '656.59'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.66%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.39'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -5.80%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.420'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -5.30%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.09%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.995'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -7.60%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '3.347.06'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -2.46%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.208'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -3.12%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '41.06'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -3.45%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '97.534.01'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.29%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.09'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -5.75%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000253'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -6.40%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.964.10'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.41%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '8.55'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -9.91%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.350.80'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -2.14%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.77'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -2.00%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.509'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.30%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.77'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -2.08%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '505.43'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -2.45%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.35'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -4.62%  '
$ws.Range('B25').Value = 'NEARProtocol'
$ws.Range('C25').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '7.03'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +8.26%  '
$ws.Range('B26').Value = 'PEPE'
$ws.Range('C26').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000199'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -4.93%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '95.90'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -5.33%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '12.25'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -7.34%  '
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '11.31'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -3.16%  '
$ws.Range('B30').Value = 'Dai'
$ws.Range('C30').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.05%  '
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.141'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -8.72%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.187'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -7.22%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.58'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +8.65%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.48%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.557'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -5.79%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '28.20'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -6.43%  '
$ws.Range('B37').Value = 'Fetch.AI'
$ws.Range('C37').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.55'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +4.84%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '7.97'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.16%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '532.38'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.90%  '
$ws.Range('B40').Value = 'USDe'
$ws.Range('C40').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.01%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.150'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -4.61%  '
$ws.Range('B42').Value = 'WhiteBITCoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '24.41'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -1.42%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.844'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -4.73%  '
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0427'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.77%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.72'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +4.34%  '
$ws.Range('B46').Value = 'MantraDAO'
$ws.Range('C46').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.62'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -1.97%  '
$ws.Range('B47').Value = 'Stacks'
$ws.Range('C47').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.22'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +5.13%  '
$ws.Range('B48').Value = 'Filecoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.53'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.43%  '
$ws.Range('B49').Value = 'Cosmos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.42'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -5.32%  '
$ws.Range('B50').Value = 'OKB'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '53.47'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +4.20%  '
$ws.Range('B51').Value = 'dogwifhat'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.13'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -8.86%  '
